# Handles float input without breaking stuff
#
# This marksheet previously failed to record the student's chosen answer
# (columns A/D on the answer-key rows) whenever the underlying score logic
# didn't deal with non-integer ("float") input cleanly, and the summary
# row (10-12) ended up with placeholder/zeroed figures ("Absent").
# This script re-populates the correct summary numbers and the per-question
# "Student Ans" columns, and removes the now-unused extra Q&A block
# (columns D/E where blank, and the whole G/H block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# ---- Summary block (rows 10-12) ------------------------------------------
# Row 10: No. right / wrong / not-attempted / max
$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 28

# Row 11: marking scheme (+ per right, - per wrong)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: totals
$ws.Range("B12").Value = 72
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "71/112"

# "No.", "Marking", "Total" row labels now pick up the header style (s=4)
$ws.Range("A10").Style = $ws.Range("A9").Style
$ws.Range("A11").Style = $ws.Range("A9").Style
$ws.Range("A12").Style = $ws.Range("A9").Style

# ---- Drop the unused third Q&A block (columns G:H), rows 15-40 -----------
$ws.Range("G15:H40").Clear()

# ---- Student-answer column (A) for the first Q&A block, rows 16-40 ------
# Style 5 = correctStyle-derived (green), 6 = incorrectStyle-derived (red),
# 7 = normalStyle-derived (black, no text) -- matches the "Correct Ans" /
# student-matched colouring already used in column B/E.
$studentAnsA = @{
  16 = "Option B"
  17 = "Option D"
  18 = "Option B"
  19 = "Option C"
  21 = "Option C"
  22 = "Option D"
  25 = "Option A"
  26 = "Option C"
  27 = "Option A"
  31 = "Option D"
  32 = "Option C"
  33 = "Option D"
  36 = "Option A"
  38 = "Option A"
  39 = "Option D"
  40 = "Option D"
}
foreach ($row in $studentAnsA.Keys) {
  $cell = $ws.Cells.Item($row, 1)
  $cell.Value = $studentAnsA[$row]
  if ($row -eq 16) {
    $cell.Style = "incorrectStyle"
  } else {
    $cell.Style = "normalStyle"
  }
}

# ---- Student-answer column (D) for the second Q&A block, rows 16-18 ------
# Only the first three rows of this block still have a matching question;
# everything below (rows 19-40) loses its D/E pair entirely.
$studentAnsD = @{
  16 = "Option A"
  17 = "Option C"
  18 = "Option D"
}
foreach ($row in $studentAnsD.Keys) {
  $cell = $ws.Cells.Item($row, 4)
  $cell.Value = $studentAnsD[$row]
  $cell.Style = "normalStyle"
}

# Remove the D/E pair on every row of the second block past row 18 (it no
# longer corresponds to a question).
$ws.Range("D19:E40").Clear()

# Tidy up the used range / dimension to match (A5:E40, not A5:H40).
$ws.Range("F1:H40").Clear()
